$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: persistResolution / COPY_INSTEAD_OF_MOVE / dsVideoPort.c / 1650
$ws.Range("A4").Value = "persistResolution"
$ws.Range("B4").Value = "COPY_INSTEAD_OF_MOVE"
$ws.Range("C4").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/devicesettings/generic/rpc/srv/dsVideoPort.c"
$ws.Range("D4").Value = 1650

# B4:C4 pick up a distinct (but visually identical) font run, matching the
# workbook's own history of slightly different Calibri font entries.
$ws.Range("B4:C4").Font.ThemeColor = 1

# Row 5: _dsGetEDIDBytes / out_of_bounds_access / dsDisplay.c / 206
$ws.Range("A5").Value = "_dsGetEDIDBytes"
$ws.Range("B5").Value = "out_of_bounds_access"
$ws.Range("C5").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/devicesettings/generic/rpc/srv/dsDisplay.c"
$ws.Range("D5").Value = 206

$ws.Range("C5").Font.ThemeColor = 1

# D5 reverts to the workbook's plain default (unstyled) numeric cell.
$ws.Range("D5").ClearFormats()
$ws.Range("D5").Value = 206

# Move the active selection the way the authored file did.
$ws.Range("C7").Select() | Out-Null
